# "Se realiza todo lo de la poligonal"
# Adds the "Proy corr Y" / "Proy corr X" columns to the Coordenadas sheet
# (computed as Proy +/- Corr), shifting the old Coord N / Coord E columns
# to the right and refreshing their cumulative values, then leaves the
# "Parametros Pol" sheet as the active tab (as in the source workbook).

$wb = $excel.ActiveWorkbook

$wsProy   = $wb.Worksheets.Item("Proyecciones")
$wsCoord  = $wb.Worksheets.Item("Coordenadas")
$wsParam  = $wb.Worksheets.Item("Parametros Pol")

# --- Coordenadas: insert two new columns (E:F) for "Proy corr Y"/"Proy corr X" ---
# This shifts the existing "Coord N"/"Coord E" columns from E:F to G:H,
# carrying the header style (bold/centered/bordered) along with them.
$wsCoord.Range("E1:F1").EntireColumn.Insert()

# Row 1 headers
$wsCoord.Range("E1").Value = "Proy corr Y"
$wsCoord.Range("F1").Value = "Proy corr X"

# Blank-row separators (rows 2,4,6,8) get a single space, matching the
# style already used by the neighbouring A:D / G:H columns on those rows.
$wsCoord.Range("E2").Value = " "
$wsCoord.Range("F2").Value = " "
$wsCoord.Range("E4").Value = " "
$wsCoord.Range("F4").Value = " "
$wsCoord.Range("E6").Value = " "
$wsCoord.Range("F6").Value = " "
$wsCoord.Range("E8").Value = " "
$wsCoord.Range("F8").Value = " "

# Updated correction values (Corr Y / Corr X) for the three legs
$wsCoord.Range("C3").Value = 0.005
$wsCoord.Range("D3").Value = -0.005
$wsCoord.Range("C5").Value = 0.007
$wsCoord.Range("D5").Value = -0.001
$wsCoord.Range("C7").Value = 0.002
$wsCoord.Range("D7").Value = -0.007

# Proy corr Y/X = Proy Y/X + Corr Y/X, for each leg row (3, 5, 7)
$wsCoord.Range("E3").Value = -40.187
$wsCoord.Range("F3").Value = -58.304
$wsCoord.Range("E5").Value = 60.452
$wsCoord.Range("F5").Value = -25.622
$wsCoord.Range("E7").Value = -20.265
$wsCoord.Range("F7").Value = 83.926

# Row 9 (closing row) separators stay blank like the rest of A:F there
$wsCoord.Range("E9").Value = " "
$wsCoord.Range("F9").Value = " "

# Coord N / Coord E (now G:H) are cumulative: next = previous + Proy corr
$wsCoord.Range("G3").Value = 1115.933
$wsCoord.Range("H3").Value = 2161.421
$wsCoord.Range("G5").Value = 1075.746
$wsCoord.Range("H5").Value = 2103.117
$wsCoord.Range("G7").Value = 1136.198
$wsCoord.Range("H7").Value = 2077.495
$wsCoord.Range("G9").Value = 1115.933
$wsCoord.Range("H9").Value = 2161.421

# --- Make "Parametros Pol" the active sheet (matches activeTab=2 / tabSelected) ---
$wsParam.Activate()
